# Auto-generated edit script applying cell value updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 924
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G9").Value = 85
$ws.Range("F12").Value = 513
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F19").Value = 36
$ws.Range("F22").Value = 0
$ws.Range("F24").Value = 197
$ws.Range("F25").Value = 517
$ws.Range("F28").Value = 543
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 3
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 116
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = 205
$ws.Range("F37").Value = 932
$ws.Range("F40").Value = 961
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 69
$ws.Range("F43").Value = 44
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 321
$ws.Range("F7").Value = 267
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 14
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 0
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 417
$ws.Range("F4").Value = 222
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 85
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 1489
$ws.Range("F11").Value = 38626
$ws.Range("G11").Value = 85
$ws.Range("F12").Value = 184
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 6
$ws.Range("F17").Value = 8079
$ws.Range("F20").Value = 532
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 93
$ws.Range("F24").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 530
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F39").Value = 801
$ws.Range("F41").Value = 350
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 44
